# Update "想去人数" (want-to-go count) figures in column F
# for both the "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 4678
    6  = 399
    8  = 926
    10 = 1174
    12 = 651
    13 = 59
    14 = 45
    15 = 26
    16 = 281
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
